$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '''244.68'
$ws.Range("G2").Formula = '''10'
$ws.Range("D3").Formula = '''23.07'
$ws.Range("G3").Formula = '''10'
$ws.Range("D4").Formula = '''5.411'
$ws.Range("G4").Formula = '''10'
$ws.Range("G5").Formula = '''10'
$ws.Range("D6").Formula = '''3.393'
$ws.Range("G6").Formula = '''10'
$ws.Range("D7").Formula = '''0.8096'
$ws.Range("G7").Formula = '''10'
$ws.Range("D8").Formula = '''0.9291'
$ws.Range("G8").Formula = '''10'
$ws.Range("B9").Formula = 'WazirX'
$ws.Range("C9").Formula = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Formula = '''0.1430'
$ws.Range("E9").Formula = '8WazirXWRX'
$ws.Range("G9").Formula = '''10'
$ws.Range("B10").Formula = 'MandalaExchangeToken'
$ws.Range("C10").Formula = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Formula = '''0.07431'
$ws.Range("E10").Formula = '9MandalaExchangeTokenMDX'
$ws.Range("G10").Formula = '''10'
$ws.Range("B11").Formula = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Formula = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Formula = '''0.03380'
$ws.Range("E11").Formula = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G11").Formula = '''10'
$ws.Range("B12").Formula = 'BitrueCoin'
$ws.Range("C12").Formula = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Formula = '''0.03039'
$ws.Range("E12").Formula = '11BitrueCoinBTR'
$ws.Range("G12").Formula = '''10'
$ws.Range("B13").Formula = 'BitMartToken'
$ws.Range("C13").Formula = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Formula = '''0.09343'
$ws.Range("E13").Formula = '12BitMartTokenBMX'
$ws.Range("G13").Formula = '''10'
$ws.Range("B14").Formula = 'MCDex'
$ws.Range("C14").Formula = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").Formula = '''3.939'
$ws.Range("E14").Formula = '13MCDexMCB'
$ws.Range("G14").Formula = '''10'
$ws.Range("B15").Formula = 'BitForexToken'
$ws.Range("C15").Formula = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Formula = '''0.001590'
$ws.Range("E15").Formula = '14BitForexTokenBF'
$ws.Range("G15").Formula = '''10'
$ws.Range("B16").Formula = 'CoinExToken'
$ws.Range("C16").Formula = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Formula = '''0.04812'
$ws.Range("E16").Formula = '15CoinExTokenCET'
$ws.Range("G16").Formula = '''10'
$ws.Range("B17").Formula = 'TigerCash'
$ws.Range("C17").Formula = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Formula = '''0.005462'
$ws.Range("E17").Formula = '16TigerCashTCH'
$ws.Range("G17").Formula = '''10'
$ws.Range("B18").Formula = 'HotbitToken'
$ws.Range("C18").Formula = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D18").Formula = '''0.004159'
$ws.Range("E18").Formula = '17HotbitTokenHTB'
$ws.Range("G18").Formula = '''10'
$ws.Range("B19").Formula = 'BitKan'
$ws.Range("C19").Formula = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").Formula = '''0.0009860'
$ws.Range("E19").Formula = '18BitKanKAN'
$ws.Range("G19").Formula = '''10'
$ws.Range("B20").Formula = 'NitroEx'
$ws.Range("C20").Formula = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D20").Formula = '''0.00007704'
$ws.Range("E20").Formula = '19NitroExNTX'
$ws.Range("G20").Formula = '''10'
$ws.Range("B21").Formula = 'LEO'
$ws.Range("C21").Formula = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D21").Formula = '''3.659'
$ws.Range("E21").Formula = '20LEOLEO'
$ws.Range("G21").Formula = '''10'
$ws.Range("B22").Formula = 'KuCoinToken'
$ws.Range("C22").Formula = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D22").Formula = '''6.465'
$ws.Range("E22").Formula = '21KuCoinTokenKCS'
$ws.Range("G22").Formula = '''10'
$ws.Range("B23").Formula = 'BTSEToken'
$ws.Range("C23").Formula = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").Formula = '''2.186'
$ws.Range("E23").Formula = '22BTSETokenBTSE'
$ws.Range("G23").Formula = '''10'
$ws.Range("B24").Formula = 'One'
$ws.Range("C24").Formula = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").Formula = '''0.01123'
$ws.Range("E24").Formula = '23OneONEBestin24h'
$ws.Range("G24").Formula = '''10'
$ws.Range("D25").Formula = '''0.3243'
$ws.Range("G25").Formula = '''10'
$ws.Range("G26").Formula = '''10'
$ws.Range("D27").Formula = '''0.0002447'
$ws.Range("G27").Formula = '''10'
$ws.Range("G28").Formula = '''10'
$ws.Range("G29").Formula = '''10'
$ws.Range("G30").Formula = '''10'
$ws.Range("G31").Formula = '''10'
$ws.Range("G32").Formula = '''10'
$ws.Range("G33").Formula = '''10'
$ws.Range("G34").Formula = '''10'
$ws.Range("G35").Formula = '''10'
$ws.Range("G36").Formula = '''10'
$ws.Range("G37").Formula = '''10'
$ws.Range("G38").Formula = '''10'
$ws.Range("G39").Formula = '''10'
$ws.Range("D40").Formula = '''0.03941'
$ws.Range("G40").Formula = '''10'
$ws.Range("D41").Formula = '''0.006213'
$ws.Range("G41").Formula = '''10'
$ws.Range("G42").Formula = '''10'
$ws.Range("G43").Formula = '''10'
$ws.Range("D44").Formula = '''0.006804'
$ws.Range("G44").Formula = '''10'
$ws.Range("D45").Formula = '''0.00005138'
$ws.Range("G45").Formula = '''10'
$ws.Range("G46").Formula = '''10'
$ws.Range("G47").Formula = '''10'
$ws.Range("D48").Formula = '''0.8554'
$ws.Range("G48").Formula = '''10'
$ws.Range("D49").Formula = '''0.002243'
$ws.Range("E49").Formula = '48BOLOBOLO'
$ws.Range("G49").Formula = '''10'
$ws.Range("D50").Formula = '''0.00002101'
$ws.Range("G50").Formula = '''10'
$ws.Range("D51").Formula = '''0.0002001'
$ws.Range("E51").Formula = '50SpecialPowerGoldSPGWorstin24h'
$ws.Range("G51").Formula = '''10'
